$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Add the "Normal (Web)" paragraph style (styleId "NormalWeb") that
#    the edited paragraphs now reference, matching it as closely as
#    possible to the style definition added to styles.xml.
# ------------------------------------------------------------------
$style = $d.Styles.Add("NormalWeb", 1)
$style.NameLocal = "Normal (Web)"
$style.BaseStyle = "Normal"
$style.Priority = 99
$style.UnhideWhenUsed = $true

$style.Font.Name = "Times New Roman"
$style.Font.NameAscii = "Times New Roman"
$style.Font.NameFarEast = "Times New Roman"
$style.Font.NameOther = "Times New Roman"
$style.Font.Size = 12
$style.Font.SizeBi = 12

$style.ParagraphFormat.SpaceBefore = 5
$style.ParagraphFormat.SpaceBeforeAuto = $true
$style.ParagraphFormat.SpaceAfter = 5
$style.ParagraphFormat.SpaceAfterAuto = $true
$style.ParagraphFormat.LineSpacingRule = 0

# ------------------------------------------------------------------
# 2) Replace the whole story body with the edited paragraphs: merged
#    runs, updated wording, "Normal (Web)" style + explicit dark-grey
#    run colour, and the reshuffled "Research paper " / "link: " /
#    hyperlink runs around the (relocated) _GoBack bookmark.
# ------------------------------------------------------------------
$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:spacing w:before="240" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
    <w:t>Using the news articles the whole process of predicting the forecast of exchange rate will be automated and for this, the research paper has used specific news articles as a dataset which are written by professionals following strict rules, which will make the news appropriate to foreign exchange and remove all the unwanted news headlines like of sports, politics, etc.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:spacing w:before="240" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
    <w:t>As usually traders see the news headlines and predict whether the rate of a particular currency will go up or down. The same thing will be done by the model, which will see the news headlines and predict, whether the currency rate will increase or decrease and at which rate.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:spacing w:before="240" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
    <w:t>They have news headlines with their timestamp containing the date and the time in seconds as in an hour they will receive 40 headlines. They used four hundred records consisting of a sequence of two to five words as the keywords for the model.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:spacing w:before="240" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
    <w:t>For the model, the occurrence of keywords in the dataset will be counted in weights to find the most effective news headline which can predict whether the currency rate will be steady, will go up, or will be down. For calculating the currency rate for the current time period t, they will use the data from time periods t-1 and t-2, as they will give the best accuracy.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:spacing w:before="240" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
    <w:t>They have made rules, by which they are able to handle continuous attributes and do not rely on Boolean tests. They have therefore more expressive power by retaining the strength of rule classifiers using comprehensible models and relatively fast learning algorithms.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:spacing w:before="240" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
    <w:t>The dataset they have used is not mentioned in the research paper, so I will try to find a dataset or any website where by using web scraping we can get the news headlines.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:spacing w:before="240" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
    <w:t xml:space="preserve">Research paper </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:color w:val="0E101A"/>
    </w:rPr>
    <w:t xml:space="preserve">link: </w:t>
  </w:r>
  <w:hyperlink r:id="rId4" w:history="1">
    <w:r>
      <w:t>https://citeseerx.ist.psu.edu/document?repid=rep1&amp;type=pdf&amp;doi=8be9790d0c4b7bceeed174ad818c6f0b36f738ad</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"/>
"@

$d.Content.InsertXML($xml)

# InsertXML drops <w:rStyle> references inside the fragment, so the
# hyperlink run's "Hyperlink" character style has to be (re)applied
# afterwards through the object model.
$d.Hyperlinks(1).Range.Style = "Hyperlink"

Write-Output "done"
